$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text/percentage updates (unambiguous as text, so a direct
# .Value assignment keeps them as strings without touching styles).
$ws.Range("D2").Value = '61.018.52'
$ws.Range("E2").Value = '  +0.58%  '
$ws.Range("D3").Value = '2.670.29'
$ws.Range("E3").Value = '  +2.12%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("E5").Value = '  +2.97%  '
$ws.Range("E6").Value = '  +0.58%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  -0.93%  '
$ws.Range("E9").Value = '  -3.37%  '
$ws.Range("E10").Value = '  +4.53%  '
$ws.Range("E11").Value = '  +1.15%  '
$ws.Range("E12").Value = '  -0.47%  '
$ws.Range("D13").Value = '3.127.38'
$ws.Range("E13").Value = '  +1.78%  '
$ws.Range("D14").Value = '61.007.25'
$ws.Range("E14").Value = '  +0.55%  '
$ws.Range("E15").Value = '  +1.65%  '
$ws.Range("E16").Value = '  +1.35%  '
$ws.Range("D17").Value = '2.668.99'
$ws.Range("E17").Value = '  +1.73%  '
$ws.Range("E18").Value = '  +0.53%  '
$ws.Range("E19").Value = '  -0.72%  '
$ws.Range("E20").Value = '  +0.03%  '
$ws.Range("E21").Value = '  +2.15%  '
$ws.Range("E22").Value = '  +0.21%  '
$ws.Range("E23").Value = '  +1.52%  '
$ws.Range("E24").Value = '  +1.13%  '
$ws.Range("E25").Value = '  +0.68%  '
$ws.Range("E26").Value = '  +0.17%  '
$ws.Range("D27").Value = '0.0₃0856'
$ws.Range("E27").Value = '  +0.80%  '
$ws.Range("E28").Value = '  -0.83%  '
$ws.Range("E29").Value = '  -0.02%  '
$ws.Range("E30").Value = '  +2.86%  '
$ws.Range("E31").Value = '  +0.13%  '
$ws.Range("E32").Value = '  +2.44%  '
$ws.Range("E33").Value = '  -1.27%  '
$ws.Range("E34").Value = '  +2.18%  '
$ws.Range("E35").Value = '  -0.16%  '
$ws.Range("E36").Value = '  +7.77%  '
$ws.Range("E37").Value = '  +0.42%  '
$ws.Range("E38").Value = '  +1.74%  '
$ws.Range("E39").Value = '  +4.12%  '
$ws.Range("E40").Value = '  -0.53%  '
$ws.Range("E41").Value = '  +0.60%  '
$ws.Range("E42").Value = '  +3.38%  '
$ws.Range("E43").Value = '  +0.35%  '
$ws.Range("E44").Value = '  +2.83%  '
$ws.Range("E45").Value = '  +1.01%  '
$ws.Range("E46").Value = '  +0.19%  '
$ws.Range("E47").Value = '  +2.34%  '
$ws.Range("E48").Value = '  -1.43%  '
$ws.Range("E49").Value = '  +5.31%  '
$ws.Range("E50").Value = '  +0.51%  '
$ws.Range("D51").Value = '1.999.80'
$ws.Range("E51").Value = '  +0.28%  '

# Price cells whose new text parses as a plain number (e.g. "529.74",
# "37.00"): a direct .Value assignment would silently convert these to
# numeric cells (and "37.00" would lose its trailing zero). Route them
# through a text formula -> copy -> paste-values-only round trip so the
# literal text is preserved exactly, without leaving a formula behind
# and without changing the cell style.
$ws.Range("D5").Formula = '="529.74"'
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("D6").Formula = '="155.48"'
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("D9").Formula = '="6.49"'
$ws.Range("D9").Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("D11").Formula = '="0.351"'
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("D15").Formula = '="22.08"'
$ws.Range("D15").Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("D18").Formula = '="4.78"'
$ws.Range("D18").Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("D19").Formula = '="354.81"'
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("D20").Formula = '="10.68"'
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("D21").Formula = '="6.33"'
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("D24").Formula = '="0.430"'
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("D26").Formula = '="0.998"'
$ws.Range("D26").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("D30").Formula = '="6.18"'
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("D32").Formula = '="1.62"'
$ws.Range("D32").Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("D33").Formula = '="149.90"'
$ws.Range("D33").Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("D35").Formula = '="1.19"'
$ws.Range("D35").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("D36").Formula = '="0.918"'
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("D37").Formula = '="0.889"'
$ws.Range("D37").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("D38").Formula = '="37.00"'
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("D39").Formula = '="305.43"'
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("D44").Formula = '="20.39"'
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("D46").Formula = '="0.999"'
$ws.Range("D46").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("D48").Formula = '="4.88"'
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("D49").Formula = '="19.32"'
$ws.Range("D49").Copy()
$ws.Range("D49").PasteSpecial(-4163)
$excel.CutCopyMode = $false
